$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.731.73"
$ws.Range("E2").Value = "  +1.56%  "
$ws.Range("D3").Value = "1.636.74"
$ws.Range("E3").Value = "  +1.54%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "213.51"
$ws.Range("E5").Value = "  +0.20%  "
$ws.Range("E6").Value = "  -0.12%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.491"
$ws.Range("E7").Value = "  +0.90%  "
$ws.Range("E8").Value = "  +0.71%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0620"
$ws.Range("E9").Value = "  +0.80%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.05"
$ws.Range("E10").Value = "  +3.00%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0834"
$ws.Range("E11").Value = "  +2.48%  "
$ws.Range("D12").Value = "1.865.42"
$ws.Range("E12").Value = "  +1.55%  "
$ws.Range("D13").Value = "1.635.09"
$ws.Range("E13").Value = "  +1.27%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.05"
$ws.Range("E14").Value = "  +0.29%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.525"
$ws.Range("E15").Value = "  +1.65%  "
$ws.Range("D16").Value = "26.703.76"
$ws.Range("E16").Value = "  +1.46%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.18"
$ws.Range("E18").Value = "  +0.74%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "208.72"
$ws.Range("E19").Value = "  +2.92%  "
$ws.Range("E20").Value = "  -0.05%  "
$ws.Range("E21").Value = "  +0.77%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.38"
$ws.Range("E22").Value = "  +0.61%  "
$ws.Range("E23").Value = "  +1.12%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.89"
$ws.Range("E24").Value = "  +0.04%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "146.05"
$ws.Range("E25").Value = "  +1.14%  "
$ws.Range("E26").Value = "  -0.15%  "
$ws.Range("E27").Value = "  -1.20%  "
$ws.Range("B28").Value = "Cosmos"
$ws.Range("C28").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.68"
$ws.Range("E28").Value = "  +1.47%  "
$ws.Range("B29").Value = "EthereumClassic"
$ws.Range("C29").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.38"
$ws.Range("E29").Value = "  +1.05%  "
$ws.Range("E31").Value = "  -0.47%  "
$ws.Range("E32").Value = "  +0.96%  "
$ws.Range("E33").Value = "  -0.26%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.50"
$ws.Range("E34").Value = "  +1.13%  "
$ws.Range("E35").Value = "  -0.62%  "
$ws.Range("D36").Value = "1.164.30"
$ws.Range("E36").Value = "  +0.29%  "
$ws.Range("E37").Value = "  +0.90%  "
$ws.Range("E38").Value = "  +2.25%  "
$ws.Range("E39").Value = "  -0.11%  "
$ws.Range("E40").Value = "  -0.08%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.502"
$ws.Range("E41").Value = "  +0.22%  "
$ws.Range("E42").Value = "  +1.28%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.38"
$ws.Range("E43").Value = "  +2.50%  "
$ws.Range("D44").Value = "1.774.25"
$ws.Range("E44").Value = "  +1.28%  "
$ws.Range("E45").Value = "  +0.89%  "
$ws.Range("E46").Value = "  +0.83%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "54.71"
$ws.Range("E47").Value = "  +0.46%  "
$ws.Range("E48").Value = "  +8.08%  "
$ws.Range("E49").Value = "  +0.69%  "
$ws.Range("E50").Value = "  +5.32%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.409"
$ws.Range("E51").Value = "  +0.52%  "
